$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor entry on Drivetrain row (row 4 header info row for F column "Sensor")
$ws.Range("F4").Value = "Sensor"

# Pigeon IMU sensor for the encoder row (row 6)
$ws.Range("F6").Value = "Pigeon IMU"

# Update the two previously-unassigned rows (13 & 14) to the new "Box" subsystem
$ws.Range("D13").Value = "Box"
$ws.Range("E13").Value = "Box position"
$ws.Range("F13").Value = "2 limit switch"

$ws.Range("D14").Value = "Box"
$ws.Range("E14").Value = "Intake roller"

# Add a new row 15 for the Pigeon (centrale inertielle) on the Drivetrain
$ws.Range("A15").Value = "CAN"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Pigeon"
$ws.Range("D15").Value = "Drivetrain"
$ws.Range("E15").Value = "Pigeon centrale inertielle"

# Update selection to match the saved workbook state
$ws.Range("E14").Select()
